# Applies the edits described by the commit "Updated obs table and finished proposal draft."
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
    return $ok
}

# 1. "providing complete coverage of the galactic disks and tidal features." -> "...complete coverage of the entire system."
Replace-Text "complete coverage of the galactic disks and tidal features." "complete coverage of the entire system." | Out-Null

# 2. Remove the AGN sentence fragment after "...internal shocks."
Replace-Text "internal shocks. It is likely that many of our targets host an active galactic nuclei and would therefore contaminate our measurements of star formation." "internal shocks." | Out-Null

# 3. "While AGN contamination ... identify such areas of contamination using other emission lines."
Replace-Text "While AGN contamination may be an issue for some of star formation measurements, we will be able identify such areas of contamination using other emission lines." "While contamination may be an issue for some of star formation measurements, we will be able identify such areas using other emission lines." | Out-Null

# 4. Remove the "(Spindler et al. 2018)" hyperlink citation, leaving just a period.
Replace-Text "specific star formation rate (Spindler et al. 2018)." "specific star formation rate." | Out-Null

# 5. Rewrite the seeing / exposure-time / observing-nights paragraph. This also removes the
#    "_GoBack" bookmark that used to sit inside this span (Word drops bookmarks contained
#    entirely within a replaced range), so it is re-created below at its new location.
$oldPara = ' of the LIFU. With the assumption of a seeing of 1.2``' +
    ' and a sky brightness of 21.4 (dark) in the V band, we will require an exposure of approximately 5400s per target for a S/N/pix = 10 and S/N/Ang = 21.61 at 22 mags / arcsec2. This calculation used the WEAVE exposure time calculator. Thus, the on target observing time for the 14 targets would be 2.62 (rounded up to 3) dark nights mid February 2024. If overheads are not included in the calculation, we would estimate that per target would require 8100s total time, requiring 3.94 (rounded to 4) nights of observing. Assuming Grey time (sky >= 20.8 mag/arcsec2), we find that an acceptable S/N/pix=17.82 and S/N/Ang=9.07 would be achieved in the same observing time.'

$newPara = ' of the LIFU. With the assumption of a FWHM seeing of 0.75``' +
    ' and a sky brightness of 21.4 (dark) in the V band, we will require an exposure of approximately 6000s per target for a S/N/pix = 10.12 and S/N/Ang = 19.88 to reach a depth of 22 mags / arcsec2 while not saturating in the galactic cores. We will collaborate with the WEAVE team on the best dither pattern. This calculation used the WEAVE exposure time calculator. Thus, the on target observing time for the 14 targets would be 2.9 (rounded up to 3) dark nights. The targets are best available in May 2024, but are observable throughout the semester. If overheads are not included in the calculation, we would estimate that per target would require 9000s total time, requiring 4.4 (rounded to 5) nights of observing. Assuming Grey time (sky >= 20.8 mag/arcsec2), we find that an acceptable S/N/pix=16.04 and S/N/Ang=8.17 would be achieved in the same observing time.'

Replace-Text $oldPara $newPara | Out-Null

# Re-insert the "_GoBack" bookmark right after "...best dither pattern. " and before
# "This calculation used the WEAVE exposure time calculator.", matching its new position.
$bmFind = $d.Content.Duplicate
$bmOk = $bmFind.Find.Execute("best dither pattern. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($bmOk) {
    $bmRange = $d.Range($bmFind.End, $bmFind.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
} else {
    Write-Output "WARNING: could not relocate _GoBack bookmark"
}

Write-Output "Done."
